$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.13040000000001
$ws.Range("A4").Value = -20.82219999999998
$ws.Range("A6").Value = -23.067
$ws.Range("A7").Value = -21.8475
$ws.Range("B7").Value = 5.195799999999999
$ws.Range("A8").Value = -22.33490000000002
$ws.Range("B11").Value = 5.337300000000003
$ws.Range("B12").Value = 5.362200000000001
$ws.Range("C12").Value = -10.96879999999999
$ws.Range("E12").Value = 17.33480000000001
$ws.Range("C13").Value = -12.76229999999999
$ws.Range("E13").Value = 16.5945
$ws.Range("C14").Value = -14.40409999999999
$ws.Range("B15").Value = 5.066599999999997
$ws.Range("A16").Value = -21.51679999999999
$ws.Range("C16").Value = -12.1129
$ws.Range("C19").Value = -11.98670000000001
$ws.Range("A20").Value = -22.2643
$ws.Range("B20").Value = 4.880199999999997
$ws.Range("C20").Value = -12.16859999999999
$ws.Range("A21").Value = -22.3763
$ws.Range("B21").Value = 5.047899999999998
$ws.Range("B22").Value = 10.26100000000001
$ws.Range("C22").Value = -12.6195
$ws.Range("E22").Value = 16.71390000000002
$ws.Range("B23").Value = 9.167700000000005
$ws.Range("E25").Value = 17.14890000000001
$ws.Range("A28").Value = -21.8213
$ws.Range("A29").Value = -21.69750000000001
$ws.Range("B29").Value = 5.558500000000002
$ws.Range("E29").Value = 17.19960000000001
$ws.Range("A30").Value = -21.67030000000002
$ws.Range("A32").Value = -21.36999999999999
$ws.Range("B34").Value = 9.629000000000007
$ws.Range("E34").Value = 16.965
$ws.Range("C36").Value = -12.3399
$ws.Range("A40").Value = -19.55449999999999
$ws.Range("B42").Value = 10.1598
$ws.Range("B43").Value = 5.382899999999998
$ws.Range("C43").Value = -13.0212
$ws.Range("E43").Value = 17.27580000000003
$ws.Range("B44").Value = 4.814000000000001
$ws.Range("B45").Value = 4.93
$ws.Range("A46").Value = -21.75320000000001
$ws.Range("B46").Value = 5.585900000000002
$ws.Range("C46").Value = -11.3272
$ws.Range("E48").Value = 17.4011
$ws.Range("B50").Value = 4.562399999999998
$ws.Range("C50").Value = -13.9517
$ws.Range("A51").Value = -22.33929999999999
$ws.Range("B51").Value = 4.906199999999999
$ws.Range("A52").Value = -22.0939
$ws.Range("A57").Value = -22.95170000000001
$ws.Range("B57").Value = 5.102999999999998
$ws.Range("A59").Value = -22.37630000000001
$ws.Range("E60").Value = 15.697
$ws.Range("A62").Value = -22.34070000000002
$ws.Range("B65").Value = 5.364299999999998
$ws.Range("A66").Value = -21.5339
$ws.Range("B66").Value = 4.738999999999997
$ws.Range("B67").Value = 4.885299999999998
$ws.Range("E68").Value = 17.71380000000001
$ws.Range("E70").Value = 18.02120000000002
$ws.Range("E71").Value = 17.32820000000001
$ws.Range("A73").Value = -19.90929999999998
$ws.Range("E73").Value = 17.02350000000001
$ws.Range("A74").Value = -21.93369999999998
$ws.Range("C76").Value = -12.298
$ws.Range("A77").Value = -20.16689999999998
$ws.Range("E78").Value = 16.62800000000002
$ws.Range("B79").Value = 10.05640000000001
$ws.Range("B84").Value = 5.4741
$ws.Range("B87").Value = 5.283699999999997
$ws.Range("E87").Value = 16.3112
$ws.Range("A92").Value = -21.34190000000002
$ws.Range("B92").Value = 4.546899999999998
$ws.Range("E92").Value = 18.78600000000002
$ws.Range("C95").Value = -11.78900000000001
$ws.Range("B97").Value = 6.238299999999999
$ws.Range("C97").Value = -11.32550000000001
$ws.Range("C99").Value = -12.3747
$ws.Range("A100").Value = -22.27460000000001
$ws.Range("E101").Value = 16.93830000000001
